{"js": "// Update the 25 \"three-digit divided by one-digit\" expressions in the\n// worksheet table to the new set of problems. Each expression is unique\n// in the document, so we do a targeted search-and-replace per pair\n// rather than touching unrelated text.\nconst replacements = [\n  [\"808\u00f79=\", \"134\u00f78=\"],\n  [\"777\u00f74=\", \"738\u00f77=\"],\n  [\"948\u00f74=\", \"143\u00f77=\"],\n  [\"110\u00f78=\", \"870\u00f72=\"],\n  [\"108\u00f79=\", \"173\u00f79=\"],\n  [\"876\u00f75=\", \"734\u00f73=\"],\n  [\"113\u00f77=\", \"117\u00f74=\"],\n  [\"314\u00f79=\", \"692\u00f74=\"],\n  [\"981\u00f79=\", \"377\u00f73=\"],\n  [\"551\u00f76=\", \"457\u00f79=\"],\n  [\"921\u00f79=\", \"146\u00f77=\"],\n  [\"273\u00f77=\", \"578\u00f72=\"],\n  [\"218\u00f74=\", \"845\u00f77=\"],\n  [\"552\u00f77=\", \"575\u00f78=\"],\n  [\"307\u00f74=\", \"230\u00f72=\"],\n  [\"460\u00f78=\", \"430\u00f73=\"],\n  [\"615\u00f72=\", \"657\u00f73=\"],\n  [\"931\u00f77=\", \"961\u00f79=\"],\n  [\"319\u00f76=\", \"490\u00f72=\"],\n  [\"191\u00f75=\", \"631\u00f76=\"],\n  [\"378\u00f78=\", \"659\u00f77=\"],\n  [\"124\u00f76=\", \"130\u00f72=\"],\n  [\"816\u00f75=\", \"742\u00f72=\"],\n  [\"426\u00f77=\", \"920\u00f74=\"],\n  [\"622\u00f78=\", \"520\u00f77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the 25 \"three-digit divided by one-digit\" expressions in the\n# worksheet table to the new set of problems. Each expression is unique\n# in the document, so a targeted Find/Replace per pair is sufficient and\n# leaves everything else untouched.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"808\u00f79=\", \"134\u00f78=\"),\n  @(\"777\u00f74=\", \"738\u00f77=\"),\n  @(\"948\u00f74=\", \"143\u00f77=\"),\n  @(\"110\u00f78=\", \"870\u00f72=\"),\n  @(\"108\u00f79=\", \"173\u00f79=\"),\n  @(\"876\u00f75=\", \"734\u00f73=\"),\n  @(\"113\u00f77=\", \"117\u00f74=\"),\n  @(\"314\u00f79=\", \"692\u00f74=\"),\n  @(\"981\u00f79=\", \"377\u00f73=\"),\n  @(\"551\u00f76=\", \"457\u00f79=\"),\n  @(\"921\u00f79=\", \"146\u00f77=\"),\n  @(\"273\u00f77=\", \"578\u00f72=\"),\n  @(\"218\u00f74=\", \"845\u00f77=\"),\n  @(\"552\u00f77=\", \"575\u00f78=\"),\n  @(\"307\u00f74=\", \"230\u00f72=\"),\n  @(\"460\u00f78=\", \"430\u00f73=\"),\n  @(\"615\u00f72=\", \"657\u00f73=\"),\n  @(\"931\u00f77=\", \"961\u00f79=\"),\n  @(\"319\u00f76=\", \"490\u00f72=\"),\n  @(\"191\u00f75=\", \"631\u00f76=\"),\n  @(\"378\u00f78=\", \"659\u00f77=\"),\n  @(\"124\u00f76=\", \"130\u00f72=\"),\n  @(\"816\u00f75=\", \"742\u00f72=\"),\n  @(\"426\u00f77=\", \"920\u00f74=\"),\n  @(\"622\u00f78=\", \"520\u00f77=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
